$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46's "phone" cell was previously entered as text; it now matches the
# numeric convention used by every other row in the sheet.
$ws.Cells.Item(46, 1).Value = 79174445

# New payment row 47: 79174445 / Cash / 2025-08-18T17:43:44 / 20 -> 20, no
# discount, no points redeemed.
#
# The "phone" number on the newest (last) row is still stored as text, and
# the unused amount/discount columns are blank text cells, just like the
# rest of the sheet. A leading apostrophe forces Excel to store a value as
# text instead of a number; ClearFormats() afterwards removes the
# "quote prefix" formatting that entering text that way applies, so the new
# cells keep the sheet's plain default style.
$ws.Cells.Item(47, 1).Value = "'79174445"
$ws.Cells.Item(47, 1).ClearFormats()

$ws.Cells.Item(47, 2).Value = "'"
$ws.Cells.Item(47, 2).ClearFormats()

$ws.Cells.Item(47, 3).Value = "Cash"

$ws.Cells.Item(47, 4).Value = "2025-08-18T17:43:44"

$ws.Cells.Item(47, 5).Value = 20

$ws.Cells.Item(47, 6).Value = "'"
$ws.Cells.Item(47, 6).ClearFormats()

$ws.Cells.Item(47, 7).Value = 20
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 9).Value = 0
